$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Rows 1-3: collapse to "0M"
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# Row 4: 103 -> 1153
$t.Cell(4, 1).Range.Text = "1153"

# Rows 6-12: updated statistic values
$t.Cell(6, 1).Range.Text = "0.00502"
$t.Cell(7, 1).Range.Text = "0.00039"
$t.Cell(8, 1).Range.Text = "0.00057"
$t.Cell(9, 1).Range.Text = "0.00012"
$t.Cell(10, 1).Range.Text = "0.00026"
$t.Cell(11, 1).Range.Text = "0.00502"
$t.Cell(12, 1).Range.Text = "0.05167"

# Rows 44-46 (tab-separated multi-value rows) collapse down to the
# single values that used to occupy rows 1-3
$t.Cell(44, 1).Range.Text = "99.99"
$t.Cell(45, 1).Range.Text = "0.05"
$t.Cell(46, 1).Range.Text = "708"
